$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume text cells keep their exact text representation
# (matches source formatting, which stores these as plain text, not numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.402.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.980.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +7.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +9.03%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +6.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.29"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +11.52%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.92%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +7.84%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.490.06"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +12.41%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +15.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.401.83"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.981.72"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +12.15%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +10.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.30%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.469"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.87"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0886"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +12.28%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.51"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +11.60%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.75"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +13.54%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.42%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +10.82%  "

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.38"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.11%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.78"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +10.02%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.26%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.84%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.53"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.67%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0674"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +11.94%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.94"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.15%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "RenzoRestakedETH"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.015.32"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.50%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.27"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.94%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.640"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.20%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.249.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.21%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.984"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.22%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.49%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.55"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.45%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +25.42%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0235"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +13.54%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.76"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.16%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.92"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +9.36%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0862"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.83%  "
